$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = "Norqulova Gulhayo Qilich qizi"
$ws.Range("B17").Value = "AA7992447"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "066"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").Value = "Navoiy viloyati"
$ws.Range("E17").Value = "Navoiy shahri"
$ws.Range("F17").Value = "Maktabgacha ta’lim tashkiloti musiqa rahbari"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "998934314533"
$ws.Range("G17").ClearFormats()
$ws.Range("H17").Value = "25-04-2024"

# Row 18
$ws.Range("A18").Value = "Qarshiboyeva Dildora"
$ws.Range("B18").Value = "AB0538735"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "067"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").Value = "Toshkent viloyati"
$ws.Range("E18").Value = "Oqqoʻrgʻon tumani"
$ws.Range("F18").Value = "Maktabgacha ta’lim tashkiloti tarbiyachisi"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "998937022727"
$ws.Range("G18").ClearFormats()
$ws.Range("H18").Value = "25-04-2024"
